$wb = $excel.ActiveWorkbook

# New daily data (rows 524-540) for the four data sheets. Column C holds the
# day's raw count, column D holds the trailing 7-day AVERAGE().
$sheet1_c = @{524=17;525=1;526=11;527=29;528=16;529=21;530=30;531=20;532=2;533=53;534=36;535=38;536=32;537=34;538=26;539=-1;540=13}
$sheet2_c = @{524=0;525=0;526=0;527=0;528=0;529=0;530=0;531=0;532=0;533=0;534=0;535=0;536=0;537=0;538=1;539=0;540=0}
$sheet3_c = @{524=7;525=15;526=4;527=0;528=35;529=50;530=8;531=4;532=12;533=5;534=36;535=23;536=4;537=17;538=10;539=14;540=0}
$sheet4_c = @{524=5;525=5;526=5;527=5;528=5;529=5;530=6;531=6;532=6;533=6;534=8;535=9;536=10;537=10;538=10;539=10;540=11}

$sheetNames = @("Nuovi casi", "Deceduti", "Dimessi   Guariti", "Ricoveri")
$sheetData  = @($sheet1_c, $sheet2_c, $sheet3_c, $sheet4_c)
$topLeftRow = @(517, 511, 512, 514)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $ws.Activate()

    $cvals = $sheetData[$i]
    foreach ($row in 524..540) {
        $ws.Cells.Item($row, 3).Value = $cvals[$row]
        $ws.Cells.Item($row, 4).Formula = "=AVERAGE(C" + ($row - 6) + ":C" + $row + ")"
    }

    # Carry the existing "Media 7 giorni" number format down onto the newly
    # filled D cells (matches the #,##0 style used by the rows above), without
    # disturbing the formulas we just wrote.
    $fmtSrc = $ws.Cells.Item(523, 4)
    $fmtDst = $ws.Range($ws.Cells.Item(524, 4), $ws.Cells.Item(540, 4))
    $fmtSrc.Copy()
    $fmtDst.PasteSpecial(-4122)  # xlPasteFormats

    # Scroll / select column C, mirroring the recorded view state as closely
    # as this host's window model allows.
    $ws.Columns.Item(3).Select()
    $excel.ActiveWindow.ScrollRow = $topLeftRow[$i]
}

# Remove the "Terapia" (Terapia intensiva) sheet entirely.
$wb.Worksheets.Item("Terapia").Delete()

# Restore "Ricoveri" as the active sheet/tab, matching the saved selection.
$wb.Worksheets.Item("Ricoveri").Activate()
